$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.61%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.74%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.074"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.85%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07959"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.81%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.161"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.55%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.028"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.53%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9317"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.99%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09872"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.95%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1879"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.02%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09061"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.36%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03623"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.15%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09921"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.09%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.56%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005712"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.69%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.441"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.62%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.166"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.31%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "8.24%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3370"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.06%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.16%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.064"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2193"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.14%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04588"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.04%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001240"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.89%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004764"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.33%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.88%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "10.25%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04912"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.93%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007789"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.83%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1396"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.53%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007817"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.84%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.24%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01145"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "11.08%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006242"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.99%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.41%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.13"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "45.14%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001803"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.41%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.41%"
